# edit.ps1
# Applies "Update automàtic: dades i banners [2026-02-06 14:32]"
# - Refreshes DATA_EXTRACCIO (col H) timestamps on sheet "Dades_Període" for rows 2-176
# - For the 12 header-probe rows (one per station block), also refreshes
#   HORA_CONSULTA_UTC (col I) and URL_FONT (col J) from 13:30Z -> 14:00Z
# - Mirrors the URL_FONT (col F) update on sheet "Estudi_Capçaleres" rows 25-36

$wb = $excel.ActiveWorkbook
$wsDades = $wb.Worksheets.Item(1)
$wsEstudi = $wb.Worksheets.Item(2)

# Row -> new DATA_EXTRACCIO (col H) value, for every data row (2..176)
$hTimestamps = @(
    "2|2026-02-06 14:28:31",
    "3|2026-02-06 14:28:32",
    "4|2026-02-06 14:28:32",
    "5|2026-02-06 14:28:32",
    "6|2026-02-06 14:28:32",
    "7|2026-02-06 14:28:35",
    "8|2026-02-06 14:28:36",
    "9|2026-02-06 14:28:36",
    "10|2026-02-06 14:28:36",
    "11|2026-02-06 14:28:36",
    "12|2026-02-06 14:28:39",
    "13|2026-02-06 14:28:40",
    "14|2026-02-06 14:28:40",
    "15|2026-02-06 14:28:40",
    "16|2026-02-06 14:28:40",
    "17|2026-02-06 14:28:43",
    "18|2026-02-06 14:28:44",
    "19|2026-02-06 14:28:44",
    "20|2026-02-06 14:28:44",
    "21|2026-02-06 14:28:44",
    "22|2026-02-06 14:28:47",
    "23|2026-02-06 14:28:48",
    "24|2026-02-06 14:28:48",
    "25|2026-02-06 14:28:48",
    "26|2026-02-06 14:28:48",
    "27|2026-02-06 14:28:51",
    "28|2026-02-06 14:28:52",
    "29|2026-02-06 14:28:52",
    "30|2026-02-06 14:28:52",
    "31|2026-02-06 14:28:52",
    "32|2026-02-06 14:28:55",
    "33|2026-02-06 14:28:56",
    "34|2026-02-06 14:28:56",
    "35|2026-02-06 14:28:56",
    "36|2026-02-06 14:28:56",
    "37|2026-02-06 14:28:59",
    "38|2026-02-06 14:29:00",
    "39|2026-02-06 14:29:00",
    "40|2026-02-06 14:29:00",
    "41|2026-02-06 14:29:00",
    "42|2026-02-06 14:29:02",
    "43|2026-02-06 14:29:04",
    "44|2026-02-06 14:29:04",
    "45|2026-02-06 14:29:04",
    "46|2026-02-06 14:29:04",
    "47|2026-02-06 14:29:06",
    "48|2026-02-06 14:29:08",
    "49|2026-02-06 14:29:08",
    "50|2026-02-06 14:29:08",
    "51|2026-02-06 14:29:08",
    "52|2026-02-06 14:29:10",
    "53|2026-02-06 14:29:12",
    "54|2026-02-06 14:29:12",
    "55|2026-02-06 14:29:12",
    "56|2026-02-06 14:29:12",
    "57|2026-02-06 14:29:14",
    "58|2026-02-06 14:29:15",
    "59|2026-02-06 14:29:15",
    "60|2026-02-06 14:29:15",
    "61|2026-02-06 14:29:15",
    "62|2026-02-06 14:29:18",
    "63|2026-02-06 14:29:19",
    "64|2026-02-06 14:29:19",
    "65|2026-02-06 14:29:19",
    "66|2026-02-06 14:29:19",
    "67|2026-02-06 14:29:22",
    "68|2026-02-06 14:29:23",
    "69|2026-02-06 14:29:23",
    "70|2026-02-06 14:29:23",
    "71|2026-02-06 14:29:23",
    "72|2026-02-06 14:29:26",
    "73|2026-02-06 14:29:27",
    "74|2026-02-06 14:29:27",
    "75|2026-02-06 14:29:27",
    "76|2026-02-06 14:29:27",
    "77|2026-02-06 14:29:30",
    "78|2026-02-06 14:29:31",
    "79|2026-02-06 14:29:31",
    "80|2026-02-06 14:29:31",
    "81|2026-02-06 14:29:31",
    "82|2026-02-06 14:29:34",
    "83|2026-02-06 14:29:35",
    "84|2026-02-06 14:29:35",
    "85|2026-02-06 14:29:35",
    "86|2026-02-06 14:29:35",
    "87|2026-02-06 14:29:38",
    "88|2026-02-06 14:29:40",
    "89|2026-02-06 14:29:40",
    "90|2026-02-06 14:29:40",
    "91|2026-02-06 14:29:40",
    "92|2026-02-06 14:29:42",
    "93|2026-02-06 14:29:44",
    "94|2026-02-06 14:29:44",
    "95|2026-02-06 14:29:44",
    "96|2026-02-06 14:29:44",
    "97|2026-02-06 14:29:46",
    "98|2026-02-06 14:29:47",
    "99|2026-02-06 14:29:47",
    "100|2026-02-06 14:29:47",
    "101|2026-02-06 14:29:47",
    "102|2026-02-06 14:29:50",
    "103|2026-02-06 14:29:51",
    "104|2026-02-06 14:29:51",
    "105|2026-02-06 14:29:51",
    "106|2026-02-06 14:29:51",
    "107|2026-02-06 14:29:53",
    "108|2026-02-06 14:29:55",
    "109|2026-02-06 14:29:55",
    "110|2026-02-06 14:29:55",
    "111|2026-02-06 14:29:55",
    "112|2026-02-06 14:29:58",
    "113|2026-02-06 14:29:59",
    "114|2026-02-06 14:29:59",
    "115|2026-02-06 14:29:59",
    "116|2026-02-06 14:29:59",
    "117|2026-02-06 14:30:02",
    "118|2026-02-06 14:30:04",
    "119|2026-02-06 14:30:04",
    "120|2026-02-06 14:30:04",
    "121|2026-02-06 14:30:04",
    "122|2026-02-06 14:30:06",
    "123|2026-02-06 14:30:08",
    "124|2026-02-06 14:30:08",
    "125|2026-02-06 14:30:08",
    "126|2026-02-06 14:30:08",
    "127|2026-02-06 14:30:10",
    "128|2026-02-06 14:30:12",
    "129|2026-02-06 14:30:12",
    "130|2026-02-06 14:30:12",
    "131|2026-02-06 14:30:12",
    "132|2026-02-06 14:30:14",
    "133|2026-02-06 14:30:16",
    "134|2026-02-06 14:30:16",
    "135|2026-02-06 14:30:16",
    "136|2026-02-06 14:30:16",
    "137|2026-02-06 14:30:18",
    "138|2026-02-06 14:30:19",
    "139|2026-02-06 14:30:19",
    "140|2026-02-06 14:30:19",
    "141|2026-02-06 14:30:19",
    "142|2026-02-06 14:30:22",
    "143|2026-02-06 14:30:24",
    "144|2026-02-06 14:30:24",
    "145|2026-02-06 14:30:24",
    "146|2026-02-06 14:30:24",
    "147|2026-02-06 14:30:26",
    "148|2026-02-06 14:30:27",
    "149|2026-02-06 14:30:27",
    "150|2026-02-06 14:30:27",
    "151|2026-02-06 14:30:27",
    "152|2026-02-06 14:30:30",
    "153|2026-02-06 14:30:31",
    "154|2026-02-06 14:30:31",
    "155|2026-02-06 14:30:31",
    "156|2026-02-06 14:30:31",
    "157|2026-02-06 14:30:34",
    "158|2026-02-06 14:30:35",
    "159|2026-02-06 14:30:35",
    "160|2026-02-06 14:30:35",
    "161|2026-02-06 14:30:35",
    "162|2026-02-06 14:30:37",
    "163|2026-02-06 14:30:38",
    "164|2026-02-06 14:30:38",
    "165|2026-02-06 14:30:38",
    "166|2026-02-06 14:30:38",
    "167|2026-02-06 14:30:41",
    "168|2026-02-06 14:30:42",
    "169|2026-02-06 14:30:42",
    "170|2026-02-06 14:30:42",
    "171|2026-02-06 14:30:42",
    "172|2026-02-06 14:30:45",
    "173|2026-02-06 14:30:46",
    "174|2026-02-06 14:30:46",
    "175|2026-02-06 14:30:46",
    "176|2026-02-06 14:30:46"
)

foreach ($entry in $hTimestamps) {
    $parts = $entry.Split("|")
    $row = [int]$parts[0]
    $stamp = $parts[1]
    $wsDades.Cells.Item($row, 8).Value = $stamp
}

# Row -> station code, for the 12 rows whose HORA_CONSULTA_UTC / URL_FONT move
# from the "13:30" period to the "14:00" period
$probeRows = @(
    "117|YA",
    "122|DG",
    "127|D4",
    "132|CI",
    "137|XS",
    "142|ZC",
    "147|XH",
    "152|XE",
    "157|UE",
    "162|XO",
    "167|VS",
    "172|D7"
)

foreach ($entry in $probeRows) {
    $parts = $entry.Split("|")
    $row = [int]$parts[0]
    $code = $parts[1]
    $wsDades.Cells.Item($row, 9).Value = "14:00"
    $wsDades.Cells.Item($row, 10).Value = "https://www.meteo.cat/observacions/xema/dades?codi=" + $code + "&dia=2026-02-06T14:00Z"
}

# Mirror the URL_FONT (col F) change on the "Estudi_Capçaleres" sheet
$estudiRows = @(
    "25|YA",
    "26|DG",
    "27|D4",
    "28|CI",
    "29|XS",
    "30|ZC",
    "31|XH",
    "32|XE",
    "33|UE",
    "34|XO",
    "35|VS",
    "36|D7"
)

foreach ($entry in $estudiRows) {
    $parts = $entry.Split("|")
    $row = [int]$parts[0]
    $code = $parts[1]
    $wsEstudi.Cells.Item($row, 6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=" + $code + "&dia=2026-02-06T14:00Z"
}

Write-Host "Applied DATA_EXTRACCIO / HORA_CONSULTA_UTC / URL_FONT updates"
